# --------------------------------------------------------------------
# feat: add 2022-Q4 data
#
# 1) "总计" (summary) sheet: insert a new row for "2022-Q4" above the
#    existing "2022-Q3" row, and renumber the index column (A) for the
#    rows that shift down.
# 2) Workbook: insert a brand-new "2022-Q4" worksheet before "2022-Q3"
#    (built as an exact style-clone of "2022-Q3" so borders / fonts /
#    sheetPr match), then overwrite its data with the Q4 fund holdings.
# --------------------------------------------------------------------

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------
# 1) Update "总计" summary sheet
# ---------------------------------------------------------------
$summary = $wb.Worksheets.Item("总计")

# Push the existing data rows down by inserting a fresh row 2.
$summary.Rows.Item(2).Insert()

# The inserted row borrows formatting from the row above (the header),
# so strip it from B:D and re-apply the bordered/bold look used by the
# rest of column A by cloning the style already sitting on A3.
$summary.Range("B2:D2").ClearFormats()
$summary.Range("A3").Copy()
$summary.Range("A2").PasteSpecial(-4122)

$summary.Range("A2").Value = 0
$summary.Range("B2").Value = "2022-Q4"
$summary.Range("C2").Value = 18
$summary.Range("D2").Value = 8.77

# Renumber the index column for the rows that shifted down one slot.
$summary.Range("A3").Value = 1
$summary.Range("A4").Value = 2

# ---------------------------------------------------------------
# 2) Add the "2022-Q4" worksheet (clone "2022-Q3" for styling, then
#    overwrite the content with the Q4 numbers)
# ---------------------------------------------------------------
$q3 = $wb.Worksheets.Item("2022-Q3")
$q3.Copy($q3)
$q4 = $wb.Worksheets.Item("2022-Q3 (2)")
$q4.Name = "2022-Q4"

# Q3 has 28 rows (1 header + 27 data); Q4 only needs 19 (1 header + 18
# data), so drop the extra rows from the bottom.
$q4.Rows("20:28").Delete()

# Force columns B:G to stay text (matches the source workbook, which
# stores these numeric-looking values as text) while we overwrite them.
$q4.Range("B2:G19").NumberFormat = "@"

$q4.Range("A1").Value = ""
$q4.Range("B1").Value = "基金代码"
$q4.Range("C1").Value = "基金名称"
$q4.Range("D1").Value = "基金规模"
$q4.Range("E1").Value = "股票总仓位"
$q4.Range("F1").Value = "仓位占比"
$q4.Range("G1").Value = "持有市值(亿元)"
$q4.Range("H1").Value = "仓位排名"

$rows = @(
    @(0,  "008985", "东方红启东三年持有期混合",     "122.37", "89.26", "2.94", "3.5977", 8),
    @(1,  "000729", "建信中小盘先锋股票A",           "31.93",  "91.13", "3.19", "1.0186", 8),
    @(2,  "530005", "建信优化配置混合A",             "24.05",  "88.67", "3.05", "0.7335", 9),
    @(3,  "169109", "东方红睿和三年定开混合A",       "19.98",  "98.28", "3.03", "0.6054", 8),
    @(4,  "169107", "东方红恒阳五年定期开放混合",     "20.25",  "97.18", "2.85", "0.5771", 8),
    @(5,  "000756", "建信潜力新蓝筹股票A",           "15.30",  "90.58", "3.12", "0.4774", 8),
    @(6,  "010225", "东方红启航三年持有期混合B",     "12.50",  "94.55", "3.04", "0.3800", 8),
    @(7,  "014967", "建信潜力新蓝筹股票C",           "10.68",  "90.58", "3.12", "0.3332", 8),
    @(8,  "013919", "建信中小盘先锋股票C",           "10.14",  "91.13", "3.19", "0.3235", 8),
    @(9,  "910022", "东方红启航三年持有期混合A",     "7.72",   "94.55", "3.04", "0.2347", 8),
    @(10, "910028", "东方红内需增长混合A",           "7.63",   "94.51", "3.06", "0.2335", 8),
    @(11, "012243", "东方红内需增长混合B",           "3.97",   "94.51", "3.06", "0.1215", 8),
    @(12, "015102", "东方红ESG可持续投资混合A",      "3.15",   "94.09", "3.10", "0.0976", 8),
    @(13, "012439", "东方红睿和三年定开混合C",       "0.49",   "98.28", "3.03", "0.0148", 8),
    @(14, "015103", "东方红ESG可持续投资混合C",      "0.43",   "94.09", "3.10", "0.0133", 8),
    @(15, "013273", "招商沪深300地产等权重指数C",    "7.43",   "94.13", "0.06", "0.0045", 8),
    @(16, "015436", "建信优化配置混合C",             "0.11",   "88.67", "3.05", "0.0034", 9),
    @(17, "161721", "招商沪深300地产等权重指数A",    "4.56",   "94.13", "0.06", "0.0027", 8)
)

$r = 2
foreach ($row in $rows) {
    $q4.Cells.Item($r, 1).Value = $row[0]
    $q4.Cells.Item($r, 2).Value = $row[1]
    $q4.Cells.Item($r, 3).Value = $row[2]
    $q4.Cells.Item($r, 4).Value = $row[3]
    $q4.Cells.Item($r, 5).Value = $row[4]
    $q4.Cells.Item($r, 6).Value = $row[5]
    $q4.Cells.Item($r, 7).Value = $row[6]
    $q4.Cells.Item($r, 8).Value = $row[7]
    $r = $r + 1
}

# Drop the forced text format now that every text cell has a string in
# it, so no stray number-format style lingers on the range.
$q4.Range("B2:G19").ClearFormats()

# Restore the original active sheet/selection (the source workbook had
# "总计" active) so the edit doesn't leave an unrelated view-state change.
$summary.Activate()
$summary.Range("A1").Select() | Out-Null

Write-Host "2022-Q4 sheet populated"
